$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Increment every value in column E (rows 2 through 33) by 1
for ($r = 2; $r -le 33; $r++) {
    $cell = $ws.Cells.Item($r, 5)
    $cell.Value2 = $cell.Value2 + 1
}
